$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be parsed as numbers
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '29.183.39'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.834.39'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '240.48'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Value = '0.6654'
$ws.Range("E6").Value = '  -3.24%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '0.07353'
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").Value = '0.2915'
$ws.Range("E9").Value = '  -2.86%  '
$ws.Range("D10").Value = '22.59'
$ws.Range("E10").Value = '  -2.78%  '
$ws.Range("D11").Value = '0.07686'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '1.839.38'
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").Value = '4.965'
$ws.Range("E13").Value = '  -1.98%  '
$ws.Range("D14").Value = '0.6644'
$ws.Range("E14").Value = '  -2.85%  '
$ws.Range("D15").Value = '83.53'
$ws.Range("E15").Value = '  -4.35%  '
$ws.Range("D16").Value = '6.077'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("D17").Value = '29.159.09'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '0.000008251'
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("D19").Value = '225.31'
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("D20").Value = '12.41'
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").Value = '7.120'
$ws.Range("E22").Value = '  -3.81%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").Value = '8.615'
$ws.Range("E25").Value = '  -1.80%  '
$ws.Range("D26").Value = '0.1391'
$ws.Range("E26").Value = '  -4.14%  '
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("D28").Value = '1.506'
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("D29").Value = '4.100'
$ws.Range("E29").Value = '  -4.25%  '
$ws.Range("D30").Value = '4.028'
$ws.Range("E30").Value = '  -2.92%  '
$ws.Range("D31").Value = '1.185'
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("D32").Value = '0.05279'
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("D33").Value = '1.867'
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("D34").Value = '0.7485'
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").Value = '1.127'
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("D36").Value = '2.683'
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = '1.309.76'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").Value = '0.01795'
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("D39").Value = '2.716'
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("D40").Value = '0.9163'
$ws.Range("E40").Value = '  -2.04%  '
$ws.Range("D41").Value = '5.926'
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("D42").Value = '0.08417'
$ws.Range("E42").Value = '  +15.70%  '
$ws.Range("D43").Value = '0.00000000132'
$ws.Range("E43").Value = '  +7.71%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").Value = '101.81'
$ws.Range("E45").Value = '  -3.02%  '
$ws.Range("D46").Value = '1.970.08'
$ws.Range("E46").Value = '  -0.82%  '
$ws.Range("D47").Value = '0.5162'
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("D48").Value = '1.763'
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("D49").Value = '63.12'
$ws.Range("E49").Value = '  -2.88%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.035'
$ws.Range("E50").Value = '  -5.18%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05930'
$ws.Range("E51").Value = '  -0.48%  '
